$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the now-obsolete order rows (old rows 5,7,8 i.e. ORDS18-004, 006, 007) ---
# and renumber / refresh the remaining rows so only 3 data rows remain (ORDS18-001..003).
# Strategy: delete rows 5 through 9 entirely (shifting nothing below them up, since
# they're the last rows), then overwrite what is left of rows 3 & 4 with the refreshed
# values that replaced the former row 6 / row 9 content.

$ws.Range("A5:N9").Delete()

# Row 3 (was "ORDS18-002 / Update Search" - now refreshed with the Doc Retrieval search order)
$ws.Cells.Item(3, 1).Value = 45440.041666608799
$ws.Cells.Item(3, 2).Value = "ORDS18-002"
$ws.Cells.Item(3, 3).Value = "SIPL5316"
$ws.Cells.Item(3, 4).Value = "SIPL5688"
$ws.Cells.Item(3, 5).Value = ""
$ws.Cells.Item(3, 6).Value = ""
$ws.Cells.Item(3, 7).Value = "Old Republic Diversified Services"
$ws.Cells.Item(3, 8).Value = "Title"
$ws.Cells.Item(3, 9).Value = "Search"
$ws.Cells.Item(3, 10).Value = "Doc Retrieval"
$ws.Cells.Item(3, 11).Value = "AL"
$ws.Cells.Item(3, 12).Value = "Autauga"
$ws.Cells.Item(3, 13).Value = "WIP"
$ws.Cells.Item(3, 14).Value = "Search(T1)"

# Row 4 (was "ORDS18-003 / 30 Years Search" - now refreshed with the Typing order)
$ws.Cells.Item(4, 1).Value = 45443.041666608799
$ws.Cells.Item(4, 2).Value = "ORDS18-003"
$ws.Cells.Item(4, 3).Value = ""
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = "SIPL5317"
$ws.Cells.Item(4, 6).Value = "SIPL5317"
$ws.Cells.Item(4, 7).Value = "Old Republic Diversified Services"
$ws.Cells.Item(4, 8).Value = "Title"
$ws.Cells.Item(4, 9).Value = "Typing"
$ws.Cells.Item(4, 10).Value = "Typing"
$ws.Cells.Item(4, 11).Value = "AL"
$ws.Cells.Item(4, 12).Value = "Autauga"
$ws.Cells.Item(4, 13).Value = "Typing"
$ws.Cells.Item(4, 14).Value = "Typing(T1)"

# --- Column width tweaks to match the refreshed data ---
$ws.Columns.Item(3).ColumnWidth = 19.72
$ws.Columns.Item(8).ColumnWidth = 7.39
$ws.Columns.Item(9).ColumnWidth = 16
$ws.Columns.Item(10).ColumnWidth = 18.28

# --- Selection moved by the user while reviewing the refreshed rows ---
$ws.Range("I8").Select()
